$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 5) formatting tweak: smaller font + shorter row ----
$ws.Range("A5:E5").Font.Size = 14
$ws.Rows.Item(5).RowHeight = 20

# ---- New Requirement Traceability Matrix rows (6-13) ----
$srsUrl = "https://github.com/sangeetajoshi/E-Appointment_Final/tree/master/Documentation/SRS"
$testUrl = "https://github.com/sangeetajoshi/E-Appointment_Final/tree/master/Testing/Test%20Cases"

$features = @(
    "User Login",
    "Make Appointment Prospect Student",
    "Make Appointment Current  Student",
    "View Appointment - Faculty/Student",
    "Set Up Appointment",
    "Edit Profile(Change Password)",
    "Admin- Ad/Update/Delete User",
    "Sign Out"
)

for ($i = 0; $i -lt $features.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $features[$i]
    $ws.Cells.Item($row, 3).Value = $srsUrl
    $ws.Cells.Item($row, 5).Value = $testUrl
    $ws.Rows.Item($row).RowHeight = 45
}

# ---- Update the saved view so the scroll position / active cell matches ----
$ws.Range("E13").Select() | Out-Null
